$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.761892557144165
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5004.145919138922
$ws.Range("I2").Value = 0.1591121006873033
$ws.Range("J2").Value = 0.1580498251401419
$ws.Range("K2").Value = 0.1569477563808311
$ws.Range("L2").Value = 0.1514760038114269
$ws.Range("M2").Value = 0.1387910209955977
$ws.Range("N2").Value = 0.128097967178595
$ws.Range("O2").Value = 0.1196284634314561
$ws.Range("P2").Value = 0.113923161534484
$ws.Range("Q2").Value = 0.1100390360560334
$ws.Range("R2").Value = 0.1072477364516625
$ws.Range("S2").Value = 0.1058260440133503
$ws.Range("T2").Value = 0.1056167265907091
$ws.Range("U2").Value = 0.1056167265907091
$ws.Range("V2").Value = 0.1055665025859655
$ws.Range("W2").Value = 0.1055665025859655
$ws.Range("X2").Value = 0.1055665025859655
$ws.Range("Y2").Value = 0.1055467040767821
$ws.Range("C3").Value = 0.9212172031402588
$ws.Range("E3").Value = 6610.398598347603
$ws.Range("H3").Value = 0.1659774421428239
$ws.Range("I3").Value = 0.1575833443593304
$ws.Range("J3").Value = 0.1483452988422972
$ws.Range("K3").Value = 0.1471429309653767
$ws.Range("L3").Value = 0.1464998121539182
$ws.Range("M3").Value = 0.1464998121539182
$ws.Range("N3").Value = 0.1456481618502464
$ws.Range("O3").Value = 0.14536687153685
$ws.Range("P3").Value = 0.14526248562215
$ws.Range("Q3").Value = 0.14526248562215
$ws.Range("R3").Value = 0.1452530009339852
$ws.Range("S3").Value = 0.1452530009339852
$ws.Range("T3").Value = 0.1451290852085406
$ws.Range("U3").Value = 0.1448954947982222
$ws.Range("V3").Value = 0.1448840077303946
$ws.Range("W3").Value = 0.1448840077303946
$ws.Range("X3").Value = 0.1448576724824094
$ws.Range("Y3").Value = 0.1448576724824094
$ws.Range("C4").Value = 0.8369781970977783
$ws.Range("E4").Value = 6626.218627873102
$ws.Range("I4").Value = 0.1611938656072738
$ws.Range("J4").Value = 0.1546682612804144
$ws.Range("K4").Value = 0.1546682612804144
$ws.Range("L4").Value = 0.1486563650517406
$ws.Range("M4").Value = 0.1472052572776577
$ws.Range("N4").Value = 0.1460995352558807
$ws.Range("O4").Value = 0.1460995352558807
$ws.Range("P4").Value = 0.1457508374445595
$ws.Range("Q4").Value = 0.1456574230548097
$ws.Range("R4").Value = 0.1456574230548097
$ws.Range("S4").Value = 0.1456574230548097
$ws.Range("T4").Value = 0.145559858344043
$ws.Range("U4").Value = 0.145559858344043
$ws.Range("V4").Value = 0.1454005385891808
$ws.Range("W4").Value = 0.1454005385891808
$ws.Range("X4").Value = 0.1452570449163436
$ws.Range("Y4").Value = 0.145166055124232
$ws.Range("C5").Value = 0.8124673366546631
$ws.Range("E5").Value = 6550.770932176638
$ws.Range("I5").Value = 0.1489576926698324
$ws.Range("J5").Value = 0.1456591353442539
$ws.Range("K5").Value = 0.1448654676888615
$ws.Range("L5").Value = 0.1448654676888615
$ws.Range("M5").Value = 0.1440916580076643
$ws.Range("N5").Value = 0.1440916580076643
$ws.Range("O5").Value = 0.1440916580076643
$ws.Range("P5").Value = 0.1440916580076643
$ws.Range("Q5").Value = 0.1440916580076643
$ws.Range("R5").Value = 0.1440916580076643
$ws.Range("S5").Value = 0.1440916580076643
$ws.Range("T5").Value = 0.143781316003586
$ws.Range("U5").Value = 0.143781316003586
$ws.Range("V5").Value = 0.143781316003586
$ws.Range("W5").Value = 0.143781316003586
$ws.Range("X5").Value = 0.143781316003586
$ws.Range("Y5").Value = 0.1436953398085115
$ws.Range("C6").Value = 0.8125030994415283
$ws.Range("E6").Value = 6611.003830806123
$ws.Range("I6").Value = 0.1590916124052761
$ws.Range("J6").Value = 0.1497386987777205
$ws.Range("K6").Value = 0.1497386987777205
$ws.Range("L6").Value = 0.1461098103723791
$ws.Range("M6").Value = 0.1456501393930245
$ws.Range("N6").Value = 0.1453273294858854
$ws.Range("O6").Value = 0.1450487739632993
$ws.Range("P6").Value = 0.1450487739632993
$ws.Range("Q6").Value = 0.1450487739632993
$ws.Range("R6").Value = 0.1450487739632993
$ws.Range("S6").Value = 0.1449515017530635
$ws.Range("T6").Value = 0.1449515017530635
$ws.Range("U6").Value = 0.1449070437887896
$ws.Range("V6").Value = 0.1449070437887896
$ws.Range("W6").Value = 0.1448694703860842
$ws.Range("X6").Value = 0.1448694703860842
$ws.Range("Y6").Value = 0.1448694703860842
$ws.Range("C7").Value = 0.7969005107879639
$ws.Range("E7").Value = 6597.564719639326
$ws.Range("I7").Value = 0.154725508759533
$ws.Range("J7").Value = 0.1466471724221224
$ws.Range("K7").Value = 0.1455559557156018
$ws.Range("L7").Value = 0.1455559557156018
$ws.Range("M7").Value = 0.1448365411940499
$ws.Range("N7").Value = 0.1448365411940499
$ws.Range("O7").Value = 0.144743445378222
$ws.Range("P7").Value = 0.144743445378222
$ws.Range("Q7").Value = 0.144743445378222
$ws.Range("R7").Value = 0.144743445378222
$ws.Range("S7").Value = 0.144743445378222
$ws.Range("T7").Value = 0.144607499408174
$ws.Range("U7").Value = 0.144607499408174
$ws.Range("V7").Value = 0.144607499408174
$ws.Range("W7").Value = 0.144607499408174
$ws.Range("X7").Value = 0.144607499408174
$ws.Range("Y7").Value = 0.144607499408174
$ws.Range("C8").Value = 0.765596866607666
$ws.Range("E8").Value = 6604.524119208962
$ws.Range("I8").Value = 0.1523912600582533
$ws.Range("J8").Value = 0.1491476697685886
$ws.Range("K8").Value = 0.145268842291306
$ws.Range("L8").Value = 0.145268842291306
$ws.Range("M8").Value = 0.145268842291306
$ws.Range("N8").Value = 0.1452438501289801
$ws.Range("O8").Value = 0.1447563880847753
$ws.Range("P8").Value = 0.1447443380672674
$ws.Range("Q8").Value = 0.1447432650976744
$ws.Range("R8").Value = 0.1447431695572533
$ws.Range("S8").Value = 0.1447431610500485
$ws.Range("T8").Value = 0.1447431602925415
$ws.Range("U8").Value = 0.1447431602250908
$ws.Range("V8").Value = 0.1447431602190848
$ws.Range("W8").Value = 0.14474316021855
$ws.Range("X8").Value = 0.1447431602185024
$ws.Range("Y8").Value = 0.1447431602184983
$ws.Range("C9").Value = 0.7656512260437012
$ws.Range("E9").Value = 6654.984178723152
$ws.Range("I9").Value = 0.1620894090517024
$ws.Range("J9").Value = 0.1491072555480821
$ws.Range("K9").Value = 0.1491072555480821
$ws.Range("L9").Value = 0.1472069467639355
$ws.Range("M9").Value = 0.1466243795488526
$ws.Range("N9").Value = 0.1461515427319236
$ws.Range("O9").Value = 0.1460152903844381
$ws.Range("P9").Value = 0.1458778275198692
$ws.Range("Q9").Value = 0.1458552990805375
$ws.Range("R9").Value = 0.145829036414188
$ws.Range("S9").Value = 0.145829036414188
$ws.Range("T9").Value = 0.145829036414188
$ws.Range("U9").Value = 0.145808795063287
$ws.Range("V9").Value = 0.145808795063287
$ws.Range("W9").Value = 0.145808795063287
$ws.Range("X9").Value = 0.1457711100432546
$ws.Range("Y9").Value = 0.1457267871096131
$ws.Range("C10").Value = 0.7968888282775879
$ws.Range("E10").Value = 6571.994020096618
$ws.Range("I10").Value = 0.1556372812781492
$ws.Range("J10").Value = 0.1456004890275751
$ws.Range("K10").Value = 0.1448479815497384
$ws.Range("L10").Value = 0.1448479815497384
$ws.Range("M10").Value = 0.1448479815497384
$ws.Range("N10").Value = 0.1445595037316112
$ws.Range("O10").Value = 0.1445595037316112
$ws.Range("P10").Value = 0.1444099677207977
$ws.Range("Q10").Value = 0.144339319459255
$ws.Range("R10").Value = 0.1443332184503955
$ws.Range("S10").Value = 0.1443326782161372
$ws.Range("T10").Value = 0.1443326303792859
$ws.Range("U10").Value = 0.1443326261434117
$ws.Range("V10").Value = 0.1443326257683321
$ws.Range("W10").Value = 0.1442570167203685
$ws.Range("X10").Value = 0.1441090452260549
$ws.Range("Y10").Value = 0.1441090452260549
$ws.Range("C11").Value = 0.7968745231628418
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 6570.226711607804
$ws.Range("I11").Value = 0.1613351017488101
$ws.Range("J11").Value = 0.1520222739381487
$ws.Range("K11").Value = 0.1501566988623868
$ws.Range("L11").Value = 0.1474027962962519
$ws.Range("M11").Value = 0.1464796391993254
$ws.Range("N11").Value = 0.1443758838636448
$ws.Range("O11").Value = 0.1443758838636448
$ws.Range("P11").Value = 0.1443758838636448
$ws.Range("Q11").Value = 0.1443158444294222
$ws.Range("R11").Value = 0.1442033363479061
$ws.Range("S11").Value = 0.1442033363479061
$ws.Range("T11").Value = 0.1442033363479061
$ws.Range("U11").Value = 0.1440745947681833
$ws.Range("V11").Value = 0.1440745947681833
$ws.Range("W11").Value = 0.1440745947681833
$ws.Range("X11").Value = 0.1440745947681833
$ws.Range("Y11").Value = 0.1440745947681833
